$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the Terrain/RGB table with the new terrain list and colors.
# Row 1 (headers) stay the same: Terrain | RGB
$ws.Range("A2").Value = "plains"
$ws.Range("B2").Value = "(79, 211, 74)"

$ws.Range("A3").Value = "farmlands"
$ws.Range("B3").Value = "(200, 189, 20)"

$ws.Range("A4").Value = "hills"
$ws.Range("B4").Value = "(155, 88, 34)"

$ws.Range("A5").Value = "mountains"
$ws.Range("B5").Value = "(183, 183, 183)"

$ws.Range("A6").Value = "forest"
$ws.Range("B6").Value = "(3, 76, 0)"

$ws.Range("A7").Value = "wetlands"
$ws.Range("B7").Value = "(74, 211, 151)"

$ws.Range("A2:B7").Select()
